$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.225.32"
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = "'1.601.38"
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'211.81"
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D7').Value = "'0.483"
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'0.249"
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('D9').Value = "'0.0613"
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').Value = "'18.15"
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('D12').Value = "'1.824.99"
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = "'1.599.90"
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('D14').Value = "'4.02"
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').Value = "'0.512"
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = "'26.201.21"
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = "'61.18"
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = "'0.0₃0728"
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = "'202.83"
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').Value = "'9.25"
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').Value = "'6.00"
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('D24').Value = "'1.94"
$ws.Range('E24').Value = '  +11.63%  '
$ws.Range('D25').Value = "'144.52"
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E27').Value = '  -8.02%  '
$ws.Range('D28').Value = "'15.15"
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('D29').Value = "'6.52"
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').Value = "'0.0490"
$ws.Range('E30').Value = '  +3.02%  '
$ws.Range('D31').Value = "'1.16"
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').Value = "'2.91"
$ws.Range('E33').Value = '  -4.46%  '
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('D35').Value = "'1.48"
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('D36').Value = "'1.143.47"
$ws.Range('E36').Value = '  +3.06%  '
$ws.Range('E37').Value = '  +5.92%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('D40').Value = "'0.784"
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').Value = "'0.495"
$ws.Range('E41').Value = '  -1.78%  '
$ws.Range('D42').Value = "'0.779"
$ws.Range('E42').Value = '  -1.82%  '
$ws.Range('D43').Value = "'5.20"
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').Value = "'1.739.59"
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').Value = "'91.76"
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('E46').Value = '  -3.19%  '
$ws.Range('D47').Value = "'54.10"
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').Value = "'0.0506"
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('D49').Value = "'0.0₇0970"
$ws.Range('E49').Value = '  -9.43%  '
$ws.Range('D50').Value = "'0.407"
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('E51').Value = '  -0.14%  '
